# Refresh the crypto "symbol list" price snapshot (GitHub Actions scrape):
# updates prices in column D, a couple of "Best/Worst in 24h" suffix tags
# in column E, and re-ranks three coins (rows 41-43) whose relative order
# changed, dragging their Coin/Link/Price/Rank cells along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text; force text type
# by setting a Text number format before assignment, then reset the style
# index back to Normal/default so no stray style attribute is introduced.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "245.83"
Set-TextValue 3 4 "24.18"
Set-TextValue 4 4 "5.284"
Set-TextValue 5 4 "0.05777"
Set-TextValue 6 4 "6.509"
Set-TextValue 7 4 "3.142"
Set-TextValue 8 4 "0.8113"
Set-TextValue 9 4 "0.8610"
Set-TextValue 11 4 "0.06964"
Set-TextValue 12 4 "0.03131"
Set-TextValue 13 4 "0.02911"
Set-TextValue 14 4 "0.09389"
Set-TextValue 15 4 "3.752"
Set-TextValue 16 4 "0.001528"
Set-TextValue 17 4 "0.04683"
Set-TextValue 18 4 "0.0005973"
$ws.Cells.Item(18, 5).Value = "17OneONEWorstin24h"
Set-TextValue 19 4 "0.006192"
Set-TextValue 20 4 "0.001237"
Set-TextValue 21 4 "0.004643"
$ws.Cells.Item(21, 5).Value = "20HotbitTokenHTB"
Set-TextValue 22 4 "0.00006103"
Set-TextValue 23 4 "3.503"
Set-TextValue 24 4 "2.148"
Set-TextValue 25 4 "0.3195"
Set-TextValue 28 4 "0.0002333"
Set-TextValue 40 4 "0.03707"
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.003044"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1055"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 43 4 "0.002792"
$ws.Cells.Item(43, 5).Value = "42CEJICEJIBestin24h"
Set-TextValue 44 4 "0.007765"
Set-TextValue 45 4 "0.00005256"
Set-TextValue 46 4 "0.00000000750"
Set-TextValue 47 4 "0.4002"
Set-TextValue 48 4 "0.002443"
Set-TextValue 49 4 "0.00002101"
Set-TextValue 50 4 "0.0002001"
